# Adds the new NFL "games" rows (weeks 11-13 of the 2023 season) that were
# uploaded in the source commit, extends the shared H-column formula
# (home_team = team2) down through the new rows, and updates the sheet's
# used range / selection bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# row, week, date(serial), team1, team2, score1, score2
# blank week/date/team/score fields mean "leave that cell empty"
$data = @"
424,11,45246,CIN,BAL,20,34
425,11,45249,LAC,GB,20,23
426,11,45249,ARI,HOU,16,21
427,11,45249,NYG,WAS,31,19
428,11,45249,DAL,CAR,33,10
429,11,45249,TEN,JAX,14,34
430,11,45249,PIT,CLE,10,13
431,11,45249,LV,MIA,13,20
432,11,45249,CHI,DET,26,31
433,11,45249,TB,SF,14,27
434,11,45249,SEA,LA,16,17
435,11,45249,NYJ,BUF,6,32
436,11,45249,MIN,DEN,20,21
437,11,45249,PHI,KC,21,17
438,12,45253,GB,DET,29,22
439,12,45253,WAS,DAL,10,45
440,12,45253,SF,SEA,31,13
441,12,45254,MIA,NYJ,34,13
442,12,45256,CAR,TEN,10,17
443,12,45256,TB,IND,20,27
444,12,45256,NO,ATL,20,27
445,12,45256,PIT,CIN,16,10
446,12,45256,NE,NYG,7,10
447,12,45256,CLE,DEN,12,29
448,12,45256,LA,ARI,37,14
449,12,45256,KC,LV,31,17
450,12,45256,BUF,PHI,34,37
451,12,45256,BAL,LAC,20,10
452,12,45256,JAX,HOU,24,21
453,12,45257,CHI,MIN,12,10
454,13,45260,,,,
455,13,45263,,,,
456,13,45263,,,,
457,13,45263,,,,
458,13,45263,,,,
459,13,45263,,,,
460,13,45263,,,,
461,13,45263,,,,
462,13,45263,,,,
463,13,45263,,,,
464,13,45263,,,,
465,13,45263,,,,
466,13,45263,,,,
467,,,,,,
468,,,,,,
469,,,,,,
"@

$lines = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

foreach ($line in $lines) {
    $f = $line.Split(",")
    $r = [int]$f[0].Trim()

    # Column A: every one of these rows is the 2023 season
    $ws.Range("A$r").Value = 2023

    if ($f[1].Trim() -ne "") { $ws.Range("B$r").Value = [int]$f[1].Trim() }
    if ($f[2].Trim() -ne "") { $ws.Range("C$r").Value = [int]$f[2].Trim() }
    if ($f[3].Trim() -ne "") { $ws.Range("D$r").Value = $f[3].Trim() }
    if ($f[4].Trim() -ne "") { $ws.Range("E$r").Value = $f[4].Trim() }
    if ($f[5].Trim() -ne "") { $ws.Range("F$r").Value = [int]$f[5].Trim() }
    if ($f[6].Trim() -ne "") { $ws.Range("G$r").Value = [int]$f[6].Trim() }

    # H = home_team, the shared formula "=E<row>" already used for every
    # other row in the table
    if ($f[3].Trim() -ne "" -and $f[4].Trim() -ne "") {
        $ws.Range("H$r").Formula = "=E$r"
    }
}

# Column C holds dates - copy the existing short-date style (already used
# all the way down the column) onto every new date cell in one shot instead
# of assigning a NumberFormat string (which would otherwise create a brand
# new, redundant cell style).
$ws.Range("C419").Copy()
$ws.Range("C424:C466").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wb.Application.Calculate()

# Move the view down to the newly-added rows, matching the author's final
# on-screen selection
$ws.Activate()
$ws.Range("D454").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 440
$win.ScrollColumn = 1
